$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in the "Recorded By" (column G) list whose two comma-separated entries
# need to be swapped (the entry that used to come second now comes first).
$rows = @(3, 6, 7, 10, 11, 12, 13, 14, 15, 17, 18, 19, 20, 21, 22, 24, 29, 32, 33, 36, 37, 38, 39, 40, 41, 43, 44, 45, 46, 47, 48, 50, 55, 58, 59, 62, 63, 64, 65, 66, 67, 69, 70, 71, 72, 73, 74, 76, 83, 84, 85, 86, 87, 90, 92, 93, 94, 96, 99, 101, 109, 110, 111, 112, 113, 116, 118, 119, 120, 122, 125, 127, 135, 136, 137, 138, 139, 142, 144, 145, 146, 148, 151, 153)

$changed = 0
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $current = [string]$cell.Value2
    $parts = $current -split ", ", 2
    if ($parts.Count -eq 2) {
        $cell.Value = "{0}, {1}" -f $parts[1], $parts[0]
        $changed++
    }
}

Write-Output ("Swapped Recorded-By order in {0} cells" -f $changed)
